$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts rows 4..66 down to 5..67
# and extends the used range dimension to A1:T67 automatically.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with a fresh weekly record, which
# duplicates the row 2 record pattern but with a new date (D) value.
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44649
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100104
$ws.Cells.Item(4, 8).Value = "Frutos de pepita"
$ws.Cells.Item(4, 9).Value = 100104003
$ws.Cells.Item(4, 10).Value = "Membrillo"
$ws.Cells.Item(4, 11).Value = "Champion"
$ws.Cells.Item(4, 12).Value = "Especial"
$ws.Cells.Item(4, 13).Value = 8
$ws.Cells.Item(4, 14).Value = 280000
$ws.Cells.Item(4, 15).Value = 280000
$ws.Cells.Item(4, 16).Value = 280000
$ws.Cells.Item(4, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(4, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(4, 19).Value = 622
$ws.Cells.Item(4, 20).Value = 450
